{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Target edit (per the supplied OOXML diff):\n//   1. \"Kickoff yard team A\" / \"Kickoff yard team B\" paragraphs: the run\n//      text is split into \"Kickoff\" + \" yard team A\"/\"B\" (wrapped with\n//      proofing (spell-check) marks) - a cosmetic run split, same visible\n//      text.\n//   2. \"Scores difference so far A-B (Indicator variable\" + \")\" (split\n//      across two runs around a \"_GoBack\" bookmark) is merged into a\n//      single run reading \"Scores difference so far A-B (Indicator\n//      variable)\"; the \"_GoBack\" bookmark is removed from this paragraph.\n//   3. Four new bulleted list items are added right after that paragraph:\n//        \"HT_fumble\"\n//        \"VT_fumble\"\n//        \"HT_rushyards/rushattempts\"\n//        \"VT_ rushyards/rushattempts\"\n//      (the \"_GoBack\" bookmark re-appears inside the last of these).\n//\n// We use Range.insertOoxml with raw WordprocessingML fragments so the\n// resulting markup (including the w:proofErr / w:bookmarkStart / w:bookmarkEnd\n// elements) matches the target precisely, rather than relying on\n// insertText/insertParagraph (which cannot produce proofErr marks or\n// control run-splitting).\n\nconst PKG_OPEN =\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>';\nconst PKG_CLOSE = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n\nfunction wrapPkg(bodyFragmentXml) {\n  return PKG_OPEN + bodyFragmentXml + PKG_CLOSE;\n}\n\n// Common pPr for the numId=2 bulleted \"Features\" list items used throughout\n// this section of the document.\nconst LIST_PPR =\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr><w:rPr><w:sz w:val=\"24\"/></w:rPr></w:pPr>';\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the three paragraphs we need to touch by their current text, so the\n// script is resilient to the exact paragraph index.\nlet idxKickoffA = -1;\nlet idxKickoffB = -1;\nlet idxScores = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t === \"Kickoff yard team A\") idxKickoffA = i;\n  else if (t === \"Kickoff yard team B\") idxKickoffB = i;\n  else if (t.indexOf(\"Scores difference so far A-B (Indicator variable\") === 0) idxScores = i;\n}\n\nif (idxKickoffA === -1 || idxKickoffB === -1 || idxScores === -1) {\n  throw new Error(\"Could not locate the target 'Features' list paragraphs.\");\n}\n\n// 1. \"Kickoff yard team A\" -> split run with proofErr marks (same text).\nconst kickoffAXml =\n  \"<w:p>\" +\n  LIST_PPR +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>Kickoff</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> yard team A</w:t></w:r>' +\n  \"</w:p>\";\nparagraphs.items[idxKickoffA]\n  .getRange(\"Whole\")\n  .insertOoxml(wrapPkg(kickoffAXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// 2. \"Kickoff yard team B\" -> split run with proofErr marks (same text).\nconst kickoffBXml =\n  \"<w:p>\" +\n  LIST_PPR +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>Kickoff</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> yard team B</w:t></w:r>' +\n  \"</w:p>\";\nparagraphs.items[idxKickoffB]\n  .getRange(\"Whole\")\n  .insertOoxml(wrapPkg(kickoffBXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// 3. \"Scores difference so far A-B (Indicator variable)\" -> single run,\n//    bookmark removed from here.\nconst scoresXml =\n  \"<w:p>\" +\n  LIST_PPR +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>Scores difference so far A-B (Indicator variable)</w:t></w:r>' +\n  \"</w:p>\";\nparagraphs.items[idxScores]\n  .getRange(\"Whole\")\n  .insertOoxml(wrapPkg(scoresXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// Re-load paragraphs since indices/anchors shifted after the replacements.\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idxScores2 = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Scores difference so far A-B (Indicator variable)\") {\n    idxScores2 = i;\n    break;\n  }\n}\nif (idxScores2 === -1) {\n  throw new Error(\"Could not re-locate the 'Scores difference...' paragraph.\");\n}\n\n// 4. Insert the four new bulleted list items right after it.\nconst newItemsXml =\n  // HT_fumble\n  \"<w:p>\" +\n  LIST_PPR +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>HT_fumble</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"</w:p>\" +\n  // VT_fumble\n  \"<w:p>\" +\n  LIST_PPR +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>VT_fumble</w:t></w:r>' +\n  \"</w:p>\" +\n  // HT_rushyards/rushattempts\n  \"<w:p>\" +\n  LIST_PPR +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>HT_rushyards</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>/</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>rushattempts</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"</w:p>\" +\n  // VT_ rushyards/rushattempts (bookmark \"_GoBack\" relocated here)\n  \"<w:p>\" +\n  LIST_PPR +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>VT_</w:t></w:r>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>rushyards</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>/</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>rushattempts</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  \"</w:p>\";\n\nparagraphs.items[idxScores2]\n  .getRange(\"Whole\")\n  .insertOoxml(wrapPkg(newItemsXml), Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# Target edit (per the supplied OOXML diff):\n#   1. \"Kickoff yard team A\" / \"Kickoff yard team B\" paragraphs: the run\n#      text is split into \"Kickoff\" + \" yard team A\"/\"B\" (wrapped with\n#      proofing (spell-check) marks) - a cosmetic run split, same visible\n#      text.\n#   2. \"Scores difference so far A-B (Indicator variable\" + \")\" (split\n#      across two runs around a \"_GoBack\" bookmark) is merged into a\n#      single run reading \"Scores difference so far A-B (Indicator\n#      variable)\"; the \"_GoBack\" bookmark is removed from this paragraph.\n#   3. Four new bulleted list items are added right after that paragraph:\n#        \"HT_fumble\"\n#        \"VT_fumble\"\n#        \"HT_rushyards/rushattempts\"\n#        \"VT_ rushyards/rushattempts\"\n#      (the \"_GoBack\" bookmark re-appears inside the last of these).\n#\n# We drive this with Range.InsertXML using raw WordprocessingML (wrapped in\n# the FlatOPC pkg:package form) so the resulting markup (including the\n# w:proofErr / w:bookmarkStart / w:bookmarkEnd elements) matches the target\n# precisely.\n\n$d = $word.ActiveDocument\n\n$pkgOpen = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>'\n$pkgClose = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n\n# Common pPr for the numId=2 bulleted \"Features\" list items used throughout\n# this section of the document.\n$listPPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"2\"/></w:numPr><w:rPr><w:sz w:val=\"24\"/></w:rPr></w:pPr>'\n\n# Locate the three paragraphs we need to touch by their current text, so the\n# script is resilient to the exact paragraph index.\n$idxKickoffA = -1\n$idxKickoffB = -1\n$idxScores = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $t = $p.Range.Text.TrimEnd()\n    if ($t -eq \"Kickoff yard team A\") { $idxKickoffA = $i }\n    elseif ($t -eq \"Kickoff yard team B\") { $idxKickoffB = $i }\n    elseif ($t.StartsWith(\"Scores difference so far A-B (Indicator variable\")) { $idxScores = $i }\n}\n\nif ($idxKickoffA -eq -1 -or $idxKickoffB -eq -1 -or $idxScores -eq -1) {\n    throw \"Could not locate the target 'Features' list paragraphs.\"\n}\n\n# 1. \"Kickoff yard team A\" -> split run with proofErr marks (same text).\n$kickoffAXml = \"<w:p>\" + $listPPr +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>Kickoff</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> yard team A</w:t></w:r>' +\n    \"</w:p>\"\n$rA = $d.Paragraphs.Item($idxKickoffA).Range\n$rA.InsertXML($pkgOpen + $kickoffAXml + $pkgClose)\n\n# 2. \"Kickoff yard team B\" -> split run with proofErr marks (same text).\n$kickoffBXml = \"<w:p>\" + $listPPr +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>Kickoff</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> yard team B</w:t></w:r>' +\n    \"</w:p>\"\n$rB = $d.Paragraphs.Item($idxKickoffB).Range\n$rB.InsertXML($pkgOpen + $kickoffBXml + $pkgClose)\n\n# 3. \"Scores difference so far A-B (Indicator variable)\" -> single run,\n#    bookmark removed from here.\n$scoresXml = \"<w:p>\" + $listPPr +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>Scores difference so far A-B (Indicator variable)</w:t></w:r>' +\n    \"</w:p>\"\n$rScores = $d.Paragraphs.Item($idxScores).Range\n$rScores.InsertXML($pkgOpen + $scoresXml + $pkgClose)\n\n# Re-locate the (now merged) \"Scores difference...\" paragraph since indices\n# may have shifted after the replacements above.\n$idxScores2 = -1\n$i = 0\nforeach ($p in $d.Paragraphs) {\n    $i = $i + 1\n    $t = $p.Range.Text.TrimEnd()\n    if ($t -eq \"Scores difference so far A-B (Indicator variable)\") { $idxScores2 = $i }\n}\nif ($idxScores2 -eq -1) {\n    throw \"Could not re-locate the 'Scores difference...' paragraph.\"\n}\n\n# 4. Insert the four new bulleted list items right after it. InsertXML at a\n#    collapsed point merges its *last* paragraph into whatever paragraph the\n#    insertion point sits in (consuming that paragraph's own content), so we\n#    first grow a disposable blank paragraph right after \"Scores...\" via\n#    InsertParagraphAfter, then target that whole (now-isolated) paragraph's\n#    Range with InsertXML \u2014 an insertXml \"whole range replace\" cleanly swaps\n#    in all of our new paragraphs with no bleed into neighbouring content.\n$scoresRange = $d.Paragraphs.Item($idxScores2).Range\n$scoresRange.Collapse(0)\n$scoresRange.InsertParagraphAfter()\n$insertionRange = $d.Paragraphs.Item($idxScores2 + 1).Range\n\n$newItemsXml =\n    # HT_fumble\n    (\"<w:p>\" + $listPPr +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>HT_fumble</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"</w:p>\") +\n    # VT_fumble\n    (\"<w:p>\" + $listPPr +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>VT_fumble</w:t></w:r>' +\n    \"</w:p>\") +\n    # HT_rushyards/rushattempts\n    (\"<w:p>\" + $listPPr +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>HT_rushyards</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>/</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>rushattempts</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"</w:p>\") +\n    # VT_ rushyards/rushattempts (bookmark \"_GoBack\" relocated here)\n    (\"<w:p>\" + $listPPr +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>VT_</w:t></w:r>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>rushyards</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>/</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:rPr><w:sz w:val=\"24\"/></w:rPr><w:t>rushattempts</w:t></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    \"</w:p>\")\n\n$insertionRange.InsertXML($pkgOpen + $newItemsXml + $pkgClose)\n"}
